$wb = $excel.ActiveWorkbook

# Update the Date value on the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2023-03-22T16:32:25+01:00"

# Clear the Condition(s) column entries on the Elements sheet for rows 4, 6 and 7
$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("AI4").Value = ""
$wsElements.Range("AI6").Value = ""
$wsElements.Range("AI7").Value = ""
